# Re-add goal-related fields that were previously dropped from the stats
# export. Insert four new rows before the current "home_xGoals" row (row 8)
# and populate them with the new metrics, pushing all subsequent rows down
# by four (old row 44 "gameresult" becomes row 48).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows starting at row 8; existing rows 8-44 shift to 12-48.
$ws.Rows("8:11").Insert()

# Copy the label formatting (bold, bordered, centered) from the row that
# used to be row 8 (now row 12, "home_xGoals") onto the newly inserted
# label cells so they match the rest of column A.
$ws.Range("A12").Copy()
$ws.Range("A8:A11").PasteSpecial(-4122)

# Populate the new rows with the goal-related fields.
$ws.Range("A8").Value = "homeGoals"
$ws.Range("B8").Value = 0
$ws.Range("A9").Value = "awayGoals"
$ws.Range("B9").Value = 0
$ws.Range("A10").Value = "homeGoalsHalfTime"
$ws.Range("B10").Value = 0
$ws.Range("A11").Value = "awayGoalsHalfTime"
$ws.Range("B11").Value = 0

$excel.CutCopyMode = $false

